$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 136.287657
$ws.Range("H2").Value = 408.862971
$ws.Range("I2").Value = 0.2628768458810872
$ws.Range("J2").Value = 0.2628768458810872
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 53.93158966666667
$ws.Range("N2").Value = 161.794769
$ws.Range("O2").Value = 0.7423286411293627
$ws.Range("P2").Value = 0.7423286411293627
$ws.Range("Q2").Value = 7350.209993955411
$ws.Range("R2").Value = 66151.8899455987
$ws.Range("S2").Value = 0.1951410117872804
$ws.Range("T2").Value = 0.1951410117872804
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 136.287657
$ws.Range("H3").Value = 408.862971
$ws.Range("I3").Value = 0.2628768458810872
$ws.Range("J3").Value = 0.2628768458810872
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.600708333333333
$ws.Range("N3").Value = 25.802125
$ws.Range("O3").Value = 0.1183824205682444
$ws.Range("P3").Value = 0.1183824205682444
$ws.Range("Q3").Value = 1172.170387290375
$ws.Range("R3").Value = 10549.53348561338
$ws.Range("S3").Value = 0.03111999732674844
$ws.Range("T3").Value = 0.03111999732674844
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 136.287657
$ws.Range("H4").Value = 408.862971
$ws.Range("I4").Value = 0.2628768458810872
$ws.Range("J4").Value = 0.2628768458810872
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.11960666666667
$ws.Range("N4").Value = 30.35882
$ws.Range("O4").Value = 0.1392889383023929
$ws.Range("P4").Value = 0.1392889383023929
$ws.Range("Q4").Value = 1379.17748236158
$ws.Range("R4").Value = 12412.59734125422
$ws.Range("S4").Value = 0.03661583676705841
$ws.Range("T4").Value = 0.03661583676705841
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 143.539174
$ws.Range("H5").Value = 430.617522
$ws.Range("I5").Value = 0.2768638492442244
$ws.Range("J5").Value = 0.2768638492442244
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 53.93158966666667
$ws.Range("N5").Value = 161.794769
$ws.Range("O5").Value = 0.7423286411293627
$ws.Range("P5").Value = 0.7423286411293627
$ws.Range("Q5").Value = 7741.295833260269
$ws.Range("R5").Value = 69671.66249934242
$ws.Range("S5").Value = 0.2055239649873098
$ws.Range("T5").Value = 0.2055239649873098
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 143.539174
$ws.Range("H6").Value = 430.617522
$ws.Range("I6").Value = 0.2768638492442244
$ws.Range("J6").Value = 0.2768638492442244
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.600708333333333
$ws.Range("N6").Value = 25.802125
$ws.Range("O6").Value = 0.1183824205682444
$ws.Range("P6").Value = 0.1183824205682444
$ws.Range("Q6").Value = 1234.538569981583
$ws.Range("R6").Value = 11110.84712983425
$ws.Range("S6").Value = 0.0327758126413728
$ws.Range("T6").Value = 0.03277581264137279
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 143.539174
$ws.Range("H7").Value = 430.617522
$ws.Range("I7").Value = 0.2768638492442244
$ws.Range("J7").Value = 0.2768638492442244
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.11960666666667
$ws.Range("N7").Value = 30.35882
$ws.Range("O7").Value = 0.1392889383023929
$ws.Range("P7").Value = 0.1392889383023929
$ws.Range("Q7").Value = 1452.559982138227
$ws.Range("R7").Value = 13073.03983924404
$ws.Range("S7").Value = 0.03856407161554179
$ws.Range("T7").Value = 0.03856407161554178
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 238.6199593333333
$ws.Range("H8").Value = 715.859878
$ws.Range("I8").Value = 0.4602593048746885
$ws.Range("J8").Value = 0.4602593048746884
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 53.93158966666667
$ws.Range("N8").Value = 161.794769
$ws.Range("O8").Value = 0.7423286411293627
$ws.Range("P8").Value = 0.7423286411293627
$ws.Range("Q8").Value = 12869.15373304202
$ws.Range("R8").Value = 115822.3835973782
$ws.Range("S8").Value = 0.3416636643547726
$ws.Range("T8").Value = 0.3416636643547725
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 238.6199593333333
$ws.Range("H9").Value = 715.859878
$ws.Range("I9").Value = 0.4602593048746885
$ws.Range("J9").Value = 0.4602593048746884
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.600708333333333
$ws.Range("N9").Value = 25.802125
$ws.Range("O9").Value = 0.1183824205682444
$ws.Range("P9").Value = 0.1183824205682444
$ws.Range("Q9").Value = 2052.300672737861
$ws.Range("R9").Value = 18470.70605464075
$ws.Range("S9").Value = 0.0544866106001232
$ws.Range("T9").Value = 0.0544866106001232
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 238.6199593333333
$ws.Range("H10").Value = 715.859878
$ws.Range("I10").Value = 0.4602593048746885
$ws.Range("J10").Value = 0.4602593048746884
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 10.11960666666667
$ws.Range("N10").Value = 30.35882
$ws.Range("O10").Value = 0.1392889383023929
$ws.Range("P10").Value = 0.1392889383023929
$ws.Range("Q10").Value = 2414.740131269329
$ws.Range("R10").Value = 21732.66118142396
$ws.Range("S10").Value = 0.06410902991979274
$ws.Range("T10").Value = 0.06410902991979273
